# Add violent crime data for 2024-05-02
# Updates the running 2024 (column K) totals, plus a handful of prior-year
# corrections (columns H and J) that ship with the same data refresh.

$wb = $excel.ActiveWorkbook

$sheetEdits = [ordered]@{
    "Citywide Totals" = @{ "K2" = 2428; "J3" = 8076; "K3" = 2329; "K4" = 489; "K5" = 152; "H6" = 7920; "K6" = 2903; "H7" = 26039; "J7" = 29282; "K7" = 8301 }
    "Grant Park" = @{ "K3" = 3; "K6" = 7 }
    "Logan Square" = @{ "K2" = 26; "K7" = 119 }
    "Austin" = @{ "K2" = 166; "K3" = 164; "K6" = 187; "K7" = 560 }
    "South Chicago" = @{ "K6" = 40; "K7" = 184 }
    "Garfield Park" = @{ "K2" = 90; "K3" = 119; "H6" = 357; "H7" = 1303; "K7" = 321 }
    "West Pullman" = @{ "K3" = 40; "K7" = 131 }
    "Grand Crossing" = @{ "K2" = 68; "K3" = 93; "K6" = 83; "K7" = 267 }
    "New City" = @{ "K2" = 57; "K7" = 196 }
    "Woodlawn" = @{ "K5" = 3; "K7" = 149 }
    "By Neighborhood" = @{ "K6" = 65; "K7" = 251; "K8" = 560; "K14" = 47; "K15" = 83; "J19" = 859; "K19" = 244; "J20" = 635; "K20" = 183; "K23" = 76; "K27" = 89; "K29" = 424; "H33" = 1303; "K33" = 321; "K36" = 101; "K37" = 267; "K38" = 7; "K40" = 20; "K42" = 290; "K48" = 101; "K51" = 87; "K52" = 230; "K53" = 119; "K54" = 156; "K55" = 89; "K59" = 13; "K63" = 29; "K65" = 196; "K67" = 328; "K72" = 39; "K73" = 81; "K77" = 60; "K78" = 110; "J79" = 800; "K83" = 184; "K84" = 59; "K85" = 404; "K87" = 6; "K89" = 110; "K90" = 71; "K91" = 78; "K94" = 100; "K95" = 131; "K96" = 111; "K99" = 149; "H101" = 26039; "J101" = 29282; "K101" = 8301 }
    "North Lawndale" = @{ "K2" = 103; "K3" = 98; "K6" = 102; "K7" = 328 }
    "South Deering" = @{ "K6" = 15; "K7" = 59 }
    "Loop" = @{ "K2" = 30; "K6" = 67; "K7" = 156 }
    "Englewood" = @{ "K2" = 111; "K3" = 141; "K6" = 138; "K7" = 424 }
    "Lake View" = @{ "K6" = 52; "K7" = 101 }
    "Chatham" = @{ "K2" = 78; "J4" = 45; "K6" = 82; "J7" = 859; "K7" = 244 }
    "Bridgeport" = @{ "K6" = 18; "K7" = 47 }
    "Ashburn" = @{ "K2" = 22; "K7" = 65 }
    "Humboldt Park" = @{ "K3" = 89; "K6" = 116; "K7" = 290 }
    "Rogers Park" = @{ "K2" = 33; "K4" = 7; "K6" = 40; "K7" = 110 }
    "Lower West Side" = @{ "K3" = 17; "K7" = 89 }
    "Douglas" = @{ "K2" = 26; "K6" = 22; "K7" = 76 }
    "West Ridge" = @{ "K4" = 3; "K7" = 111 }
    "Washington Park" = @{ "K6" = 17; "K7" = 78 }
    "Roseland" = @{ "J4" = 44; "J7" = 800 }
    "Chicago Lawn" = @{ "K2" = 63; "J3" = 208; "J7" = 635; "K7" = 183 }
    "Grand Boulevard" = @{ "K2" = 41; "K6" = 21; "K7" = 101 }
    "Auburn Gresham" = @{ "K2" = 82; "K6" = 66; "K7" = 251 }
    "West Loop" = @{ "K3" = 17; "K6" = 46; "K7" = 100 }
    "Brighton Park" = @{ "K3" = 21; "K7" = 83 }
    "Portage Park" = @{ "K4" = 4; "K7" = 81 }
    "Montclare" = @{ "K2" = 4; "K7" = 13 }
    "Uptown" = @{ "K6" = 35; "K7" = 110 }
    "Edgewater" = @{ "K2" = 24; "K6" = 36; "K7" = 89 }
    "Washington Heights" = @{ "K2" = 28; "K7" = 71 }
    "Little Italy, UIC" = @{ "K3" = 24; "K7" = 87 }
    "South Shore" = @{ "K2" = 147; "K3" = 136; "K6" = 94; "K7" = 404 }
    "Old Town" = @{ "K6" = 22; "K7" = 39 }
    "Riverdale" = @{ "K3" = 21; "K7" = 60 }
    "Hegewisch" = @{ "K3" = 9; "K7" = 20 }
    "Little Village" = @{ "K2" = 64; "K3" = 53; "K6" = 96; "K7" = 230 }
    "Ukrainian Village" = @{ "K6" = 3; "K7" = 6 }
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $sheetEdits[$sheetName]
    foreach ($ref in $cells.Keys) {
        $ws.Range($ref).Value = $cells[$ref]
    }
}

Write-Output ("Updated {0} cells across {1} sheets" -f 169, 43)
